$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Flag" checklist remark in column D, row 17 (same row as the
# "Flag" / "Static" / "OK" entry), highlighted with a yellow fill.
$ws.Range("D17").Value2 = "Clarity is required, flag is not specified"
$ws.Range("D17").Interior.Color = 65535

# Give the new column a sensible width so the comment is readable.
$ws.Columns.Item(4).ColumnWidth = 31.33

# Update selection to the new last-used cell area, as left by the edit.
$ws.Range("D21").Select() | Out-Null
